# Reorders the KPI rows (2-8) on the active sheet to match the target layout,
# updating cell values accordingly. Row contents for LPE and avg_loan_size were
# also recomputed with refreshed statistics; other KPI rows keep the same
# figures but move to new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A2').Value = 'LPE'
$ws.Range("B2").Value = 0.2097323972323972
$ws.Range("C2").Value = 0.2600773895815682
$ws.Range("D2").Value = 0.2575837742504409
$ws.Range("E2").Value = 0.2367588516193024
$ws.Range("F2").Value = 0.04785137701804368
$ws.Range("G2").Value = -0.02331853796226577
$ws.Range("H2").Value = 0.07116991498030945
$ws.Range("I2").Value = 0.07116991498030945
$ws.Range("J2").Value = 1.353594848429302
$ws.Range("K2").Value = 0.2020738920219552
$ws.Range("L2").Value = $false
$ws.Range('A3').Value = 'avg_loan_size'
$ws.Range("B3").Value = 6454.523894675925
$ws.Range("C3").Value = 7109.42636114474
$ws.Range("D3").Value = 6837.565952380953
$ws.Range("E3").Value = 6976.754165449623
$ws.Range("F3").Value = 459.5109267526456
$ws.Range("G3").Value = -132.672195695121
$ws.Range("H3").Value = 592.1831224477666
$ws.Range("I3").Value = 515.714253400145
$ws.Range("J3").Value = 0.4205325582431773
$ws.Range("K3").Value = 0.6850307076071457
$ws.Range("L3").Value = $false
$ws.Range('A4').Value = 'dq29_pot30_payment_rate_$_up_to_day'
$ws.Range("B4").Value = 0.269413401451355
$ws.Range("C4").Value = 0.2532511624585909
$ws.Range("D4").Value = 0.2877388143008039
$ws.Range("E4").Value = 0.2613587425545633
$ws.Range("F4").Value = 0.01832541284944889
$ws.Range("G4").Value = 0.008107580095972395
$ws.Range("H4").Value = 0.01021783275347649
$ws.Range("I4").Value = 0.01021783275347654
$ws.Range("J4").Value = 1.913528058943781
$ws.Range("K4").Value = 0.08687762355369738
$ws.Range("L4").Value = $false
$ws.Range('A5').Value = 'dq29_pot30_payment_rate_unit_per_day'
$ws.Range("B5").Value = 0.009300333755907537
$ws.Range("C5").Value = 0.008494525107453315
$ws.Range("D5").Value = 0.006580821059234976
$ws.Range("E5").Value = 0.007219387983925
$ws.Range("F5").Value = -0.002719512696672561
$ws.Range("G5").Value = -0.001275137123528314
$ws.Range("H5").Value = -0.001444375573144247
$ws.Range("I5").Value = -0.001444375573144247
$ws.Range("J5").Value = -1.055239397200096
$ws.Range("K5").Value = 0.3163572383656994
$ws.Range("L5").Value = $false
$ws.Range('A6').Value = 'dq29_pot30_payment_rate_unit_up_to_day'
$ws.Range("B6").Value = 0.6986772256079825
$ws.Range("C6").Value = 0.6483800661477437
$ws.Range("D6").Value = 0.3917674590808589
$ws.Range("E6").Value = 0.4579447960458904
$ws.Range("F6").Value = -0.3069097665271237
$ws.Range("G6").Value = -0.1904352701018532
$ws.Range("H6").Value = -0.1164744964252705
$ws.Range("I6").Value = -0.1164744964252704
$ws.Range("J6").Value = -4.622789547901371
$ws.Range("K6").Value = 0.001487504418300887
$ws.Range("L6").Value = $true
$ws.Range('A7').Value = 'dq30_pct_$'
$ws.Range("B7").Value = 0.6884552711616436
$ws.Range("C7").Value = 0.7361260219064817
$ws.Range("D7").Value = 0.6760494110884452
$ws.Range("E7").Value = 0.7419096882689794
$ws.Range("F7").Value = -0.01240586007319838
$ws.Range("G7").Value = 0.005783666362497746
$ws.Range("H7").Value = -0.01818952643569612
$ws.Range("I7").Value = -0.01818952643569616
$ws.Range("J7").Value = -1.402249483629397
$ws.Range("K7").Value = 0.1950165041395878
$ws.Range("L7").Value = $false
$ws.Range('A8').Value = 'dq30_pct_unit'
$ws.Range("B8").Value = 0.04917201159118303
$ws.Range("C8").Value = 0.04555056957470473
$ws.Range("D8").Value = 0.03408316209696947
$ws.Range("E8").Value = 0.03410481313023526
$ws.Range("F8").Value = -0.01508884949421356
$ws.Range("G8").Value = -0.01144575644446947
$ws.Range("H8").Value = -0.003643093049744094
$ws.Range("I8").Value = -0.003643093049744098
$ws.Range("J8").Value = -1.642464383760724
$ws.Range("K8").Value = 0.1326317532753488
$ws.Range("L8").Value = $false
